# Activity 4 slide ("Compare and contrast the self-archiving policies for
# three journals") is edited to read "...for two journals" instead of
# "...for three journals". PowerPoint splits the paragraph's single run
# into three runs at the point of the edit.

$p = $ppt.ActivePresentation

$targetSlide = $null
$targetShape = $null

foreach ($sl in $p.Slides) {
    foreach ($shp in $sl.Shapes) {
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            $t = $shp.TextFrame.TextRange.Text
            if ($t -like "*Compare and contrast the self-archiving policies for three journals*") {
                $targetSlide = $sl
                $targetShape = $shp
            }
        }
    }
}

$tr = $targetShape.TextFrame.TextRange
$para = $tr.Paragraphs(2, 1)

$oldFragment = "for three "
$newFragment = "for two "

$startPos = $para.Text.IndexOf($oldFragment) + 1
$sub = $para.Characters($startPos, $oldFragment.Length)
$sub.Text = $newFragment
